$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Question Validation Succeed")

$ws.Cells.Item(12, 1).Value = "Succeed-validationCriteria-mandatory-2"
$ws.Cells.Item(11, 1).Value = "Succeed-validationCriteria-mandatory-1"
$ws.Cells.Item(11, 11).Value = '{ "mandatory": {"encounterType":"admission"} }'
$ws.Cells.Item(12, 11).Value = '{ "mandatory": {"encounterType":["admission","surveyResponse"]} }'

$ws.Cells.Item(11, 2).Value = "SurveyAnswer"
$ws.Cells.Item(11, 3).Value = "SurveyAnswer: Full config"
$ws.Cells.Item(11, 16).Value = '{ "source": "xyz" }'

$ws.Cells.Item(12, 2).Value = "SurveyAnswer"
$ws.Cells.Item(12, 3).Value = "SurveyAnswer: Full config"
$ws.Cells.Item(12, 16).Value = '{ "source": "xyz" }'

$ws.Columns.Item(1).ColumnWidth = 56

$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

$ws2 = $wb.Worksheets.Item("Question Validation Fail")
$ws2.Range("K12").Select()

$ws.Activate()
$ws.Range("P17").Select()
